# Fix up the "基金受益憑證" (fund) sheet: give it a proper header row and
# append the seven metadata columns (property_category..index) that the
# other property sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("基金受益憑證")

# --- Row 1: proper column headers (previously this row just duplicated
# row 2's data instead of holding header labels) ---
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "dealer"
$ws.Cells.Item(1, 5).Value = "quantity"
$ws.Cells.Item(1, 6).Value = "face_value"
$ws.Cells.Item(1, 7).Value = "currency"
$ws.Cells.Item(1, 8).Value = "total"
$ws.Cells.Item(1, 9).Value = "property_category"
$ws.Cells.Item(1, 10).Value = "category"
$ws.Cells.Item(1, 11).Value = "date"
$ws.Cells.Item(1, 12).Value = "legislator_name"
$ws.Cells.Item(1, 13).Value = "legislator_id"
$ws.Cells.Item(1, 14).Value = "source_file"
$ws.Cells.Item(1, 15).Value = "index"

# --- Rows 2-9: append the metadata columns I:O, matching each row's A
# column (the per-item index number) ---
$rows = @(2, 3, 4, 5, 6, 7, 8, 9)
foreach ($r in $rows) {
    $idx = $ws.Cells.Item($r, 1).Value2

    $ws.Cells.Item($r, 9).Value = "fund"
    $ws.Cells.Item($r, 10).Value = "normal"
    # Leading apostrophe forces this to stay text instead of being
    # auto-parsed into a date serial number.
    $ws.Cells.Item($r, 11).Value = "'2012-02-01"
    $ws.Cells.Item($r, 12).Value = "吳宜臻"
    $ws.Cells.Item($r, 13).Value = 1735
    $ws.Cells.Item($r, 14).Value = "tmp2691"
    $ws.Cells.Item($r, 15).Value = $idx
}
